# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the "La Palma" / "Lanzarote" rows (city name + Muertes count)
$cityA56 = $ws.Range("A56").Value2
$cityA57 = $ws.Range("A57").Value2
$ws.Range("A56").Value2 = $cityA57
$ws.Range("A57").Value2 = $cityA56

$deathsE56 = $ws.Range("E56").Value2
$deathsE57 = $ws.Range("E57").Value2
$ws.Range("E56").Value2 = $deathsE57
$ws.Range("E57").Value2 = $deathsE56

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 6 de Abril de 2020 a las 15:22"
